$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The invoice row currently at row 4 (850778 / 5000) needs to move to the
# bottom of the list (row 7), with the rows below it (5, 6, 7) shifting up
# to take rows 4, 5, 6. A new, still-empty row 8 is appended below (only
# its amount cell carries column B's number-format style).

# Stash row 4's data out of the way first.
$ws.Range("A4:B4").Cut($ws.Range("A9:B9"))

# Shift rows 5, 6, 7 up into 4, 5, 6 (one row at a time keeps each Cut a
# simple move so cells that aren't part of the move are left untouched).
$ws.Range("A5:B5").Cut($ws.Range("A4:B4"))
$ws.Range("A6:B6").Cut($ws.Range("A5:B5"))
$ws.Range("A7:B7").Cut($ws.Range("A6:B6"))

# Drop the stashed row 4 data into the new last row (7), then clear the
# now-unused scratch row 9 completely (content + leftover formatting).
$ws.Range("A9:B9").Cut($ws.Range("A7:B7"))
$ws.Range("A9:B9").Clear()

# Add the new (still blank) row 8, carrying column B's amount-column
# number format so it matches the rest of the column.
$ws.Range("B8").NumberFormat = $ws.Range("B7").NumberFormat

# Keep the header selection in sync with the sheet's new used range.
$ws.Range("A2:B8").Select() | Out-Null
